# Auto-generated script to apply Zeromus_Profits market-price refresh values
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 182.27272
$ws.Range("I55").Value = 70.75
$ws.Range("J55").Value = 479.66666
$ws.Range("K55").Value = 70.75
$ws.Range("L55").Value = 479.66666
$ws.Range("M55").Value = 143.25
$ws.Range("N55").Value = -907.66666
$ws.Range("H116").Value = 5153474.5
$ws.Range("I116").Value = 5450682.5
$ws.Range("J116").Value = 1866.6666
$ws.Range("K116").Value = 5450682.5
$ws.Range("L116").Value = 1866.6666
$ws.Range("M116").Value = -5447240.5
$ws.Range("N116").Value = -8750.6666
$ws.Range("H128").Value = 8480
$ws.Range("J128").Value = 8480
$ws.Range("L128").Value = 8480
$ws.Range("N128").Value = -18440
$ws.Range("H133").Value = 42500
$ws.Range("J133").Value = 42500
$ws.Range("L133").Value = 42500
$ws.Range("N133").Value = -52620

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1597.7059
$ws.Range("I45").Value = 1654.3572
$ws.Range("J45").Value = 1333.3334
$ws.Range("K45").Value = 1654.3572
$ws.Range("L45").Value = 1333.3334
$ws.Range("M45").Value = -1277.3572
$ws.Range("N45").Value = -2087.3334
$ws.Range("H74").Value = 4239291
$ws.Range("I74").Value = 6946351
$ws.Range("J74").Value = 2154
$ws.Range("K74").Value = 6946351
$ws.Range("L74").Value = 2154
$ws.Range("M74").Value = -6945477
$ws.Range("N74").Value = -3902
$ws.Range("H77").Value = 4239291
$ws.Range("I77").Value = 6946351
$ws.Range("J77").Value = 2154
$ws.Range("K77").Value = 34731755
$ws.Range("L77").Value = 10770
$ws.Range("M77").Value = -34727387
$ws.Range("N77").Value = -19506
$ws.Range("H105").Value = 40370
$ws.Range("J105").Value = 40370
$ws.Range("L105").Value = 40370
$ws.Range("N105").Value = -47358
$ws.Range("H107").Value = 42000
$ws.Range("J107").Value = 42000
$ws.Range("L107").Value = 42000
$ws.Range("N107").Value = -49680
$ws.Range("H109").Value = 53333.332
$ws.Range("J109").Value = 53333.332
$ws.Range("L109").Value = 53333.332
$ws.Range("N109").Value = -56107.332
$ws.Range("H122").Value = 1753.579
$ws.Range("I122").Value = 1119.3334
$ws.Range("J122").Value = 2537.0588
$ws.Range("K122").Value = 3358.0002
$ws.Range("L122").Value = 7611.176399999999
$ws.Range("M122").Value = -908.0001999999999
$ws.Range("N122").Value = -12511.1764
$ws.Range("H123").Value = 46122.57
$ws.Range("J123").Value = 46122.57
$ws.Range("L123").Value = 46122.57
$ws.Range("N123").Value = -55922.57
$ws.Range("H125").Value = 41478.75
$ws.Range("J125").Value = 41478.75
$ws.Range("L125").Value = 41478.75
$ws.Range("N125").Value = -51318.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3351253
$ws.Range("I31").Value = 4786560
$ws.Range("J31").Value = 2202.7222
$ws.Range("K31").Value = 4786560
$ws.Range("L31").Value = 2202.7222
$ws.Range("M31").Value = -4786265
$ws.Range("N31").Value = -2792.7222
$ws.Range("H34").Value = 3351253
$ws.Range("I34").Value = 4786560
$ws.Range("J34").Value = 2202.7222
$ws.Range("K34").Value = 4786560
$ws.Range("L34").Value = 2202.7222
$ws.Range("M34").Value = -4786358
$ws.Range("N34").Value = -2606.7222
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("H124").Value = 9993.143
$ws.Range("I124").Value = 2250
$ws.Range("J124").Value = 13090.4
$ws.Range("K124").Value = 2250
$ws.Range("L124").Value = 13090.4
$ws.Range("M124").Value = 205
$ws.Range("N124").Value = -18000.4
$ws.Range("N106").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1080
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1080
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 9720
$ws.Range("N132").Value = -14780
$ws.Range("M132").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("H126").Value = 2233.3462
$ws.Range("I126").Value = 1493.9375
$ws.Range("J126").Value = 3416.4
$ws.Range("K126").Value = 4481.8125
$ws.Range("L126").Value = 10249.2
$ws.Range("M126").Value = -2011.8125
$ws.Range("N126").Value = -15189.2
$ws.Range("N105").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2256.923
$ws.Range("I7").Value = 1876.8572
$ws.Range("J7").Value = 2700.3333
$ws.Range("K7").Value = 1876.8572
$ws.Range("L7").Value = 2700.3333
$ws.Range("M7").Value = -1764.8572
$ws.Range("N7").Value = -2924.3333
$ws.Range("H122").Value = 5116.6665
$ws.Range("I122").Value = 5375
$ws.Range("J122").Value = 4600
$ws.Range("K122").Value = 16125
$ws.Range("L122").Value = 13800
$ws.Range("M122").Value = -13675
$ws.Range("N122").Value = -18700
$ws.Range("H126").Value = 2256.923
$ws.Range("I126").Value = 1876.8572
$ws.Range("J126").Value = 2700.3333
$ws.Range("K126").Value = 5630.571599999999
$ws.Range("L126").Value = 8100.999899999999
$ws.Range("M126").Value = -3160.571599999999
$ws.Range("N126").Value = -13040.9999
$ws.Range("H136").Value = 2888.754
$ws.Range("I136").Value = 3790.7
$ws.Range("J136").Value = 1445.64
$ws.Range("K136").Value = 11372.1
$ws.Range("L136").Value = 4336.92
$ws.Range("M136").Value = -8822.099999999999
$ws.Range("N136").Value = -9436.92

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 18533.334
$ws.Range("J64").Value = 18533.334
$ws.Range("L64").Value = 18533.334
$ws.Range("N64").Value = -19029.334
$ws.Range("H67").Value = 18533.334
$ws.Range("J67").Value = 18533.334
$ws.Range("L67").Value = 18533.334
$ws.Range("N67").Value = -20249.334
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
$ws.Range("H122").Value = 1457.3793
$ws.Range("I122").Value = 1060.375
$ws.Range("J122").Value = 1946
$ws.Range("K122").Value = 3181.125
$ws.Range("L122").Value = 5838
$ws.Range("M122").Value = -731.125
$ws.Range("N122").Value = -10738
$ws.Range("H123").Value = 47285.297
$ws.Range("J123").Value = 47285.297
$ws.Range("L123").Value = 47285.297
$ws.Range("N123").Value = -57085.297
$ws.Range("H126").Value = 1765.9166
$ws.Range("I126").Value = 1723.5518
$ws.Range("J126").Value = 1941.4286
$ws.Range("K126").Value = 5170.6554
$ws.Range("L126").Value = 5824.2858
$ws.Range("M126").Value = -2700.6554
$ws.Range("N126").Value = -10764.2858
